$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting rows 14:147 down to 15:148
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new data record
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C14").Value = 'Los Lagos'
$ws.Range("D14").Value = '2021-10-21'
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 100112017
$ws.Range("G14").Value = 'Apio'
$ws.Range("H14").Value = 'Americana (o)'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 35
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 11000
$ws.Range("M14").Value = 11000
$ws.Range("N14").Value = '$/docena de matas'
$ws.Range("O14").Value = 'Región de Coquimbo'
$ws.Range("P14").Value = 1833
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 'Hortaliza'
